# update figs with fontsize
# Strip the "results\" prefix (leftover from the local filesystem path used
# when the benchmark results were collected) from the model names in
# column A, and normalize any remaining backslash path separators to
# forward slashes (e.g. "results\Geotrend\bert-base-10lang-cased"
# -> "Geotrend/bert-base-10lang-cased").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '^results\\', ''
        $newVal = $newVal -replace '\\', '/'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
